# Update the "Förändrad" (Changed) date column (C) for rows 2-11
# from the old value (45224 -> 2023-10-25) to the new value (45233 -> 2023-11-03).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 3).Value = 45233
}
